$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the work-time value for the last logged day (row 10, column B)
$ws.Range("B10").Value = "13 Hours 8 Minutes"

# Move the active cell selection to D9, matching the author's last cursor position
$ws.Range("D9").Select()
